$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset (BSL NO, BRAND, ISL NO, Item Name, UOM) for rows 2..27
$data = @(
    @(4, 'Desodin', 1, 'Desodin 60ml Syrup', '60 ml'),
    @(5, 'Dinafex', 2, 'Dinafex 60mg Tablet', '30''s'),
    @(5, 'Dinafex', 3, 'Dinafex 120mg Tablet', '30''s'),
    @(5, 'Dinafex', 4, 'Dinafex 180mg Tablet', '30''s'),
    @(6, 'Dorenta', 5, 'Dorenta 50mg Tablet', '50''s'),
    @(7, 'Etorix', 6, 'Etorix 120mg Tablet', '20''s'),
    @(7, 'Etorix', 7, 'Etorix 60mg Tablet - 40''s', '40''s'),
    @(7, 'Etorix', 8, 'Etorix 90mg Tablet', '30''s'),
    @(8, 'Fenobac', 9, 'Fenobac 100ml Syrup', '100ml'),
    @(9, 'Flucloxin', 10, 'Flucloxin 500mg Capsule', '30 ''s'),
    @(9, 'Flucloxin', 11, 'Flucloxin 500mg Capsule - 36''s', '36 ''s'),
    @(10, 'Geminox', 12, 'Geminox 320mg Tablet - 8''s', '8 ''s'),
    @(11, 'Ketonic', 13, 'Ketonic 30mg IM/IV Injection - 4''s', '4''s'),
    @(11, 'Ketonic', 14, 'Ketonic 10mg Tablet', '20''s'),
    @(11, 'Ketonic', 15, 'Ketonic 30mg Injection', '5 ''s'),
    @(12, 'Kynol', 16, 'Kynol TR 200mg Capsule', '30 ''s'),
    @(12, 'Kynol', 17, 'Kynol TR 100mg Capsule', '50 ''s'),
    @(12, 'Kynol', 18, 'Kynol D 25mg Tablet', '60 ''s'),
    @(17, 'Naprox', 19, 'Naprox Plus 500mg Tablet - 30''s', '30 ''s'),
    @(19, 'Oradin', 20, 'Oradin Plus Tablet - 40''s', '40 ''s'),
    @(20, 'Osticare', 21, 'Osticare Tablet 24''s', '24''s'),
    @(23, 'Rupaday', 22, 'Rupaday Oral Solution 60ml', '1''s'),
    @(35, 'Zithrox', 23, 'Zithrox 30ml Dry Suspension', '30ml'),
    @(35, 'Zithrox', 24, 'Zithrox 500mg Tablet', '6 ''s'),
    @(35, 'Zithrox', 25, 'Zithrox 15ml Suspension', '15 ml'),
    @(35, 'Zithrox', 26, 'Zithrox 250mg Tablet - 6''s', '6''s')
)

# Clear the previous data range (old sheet had rows through 30) and any stray
# content beyond the new table so the used range / dimension shrinks correctly.
$ws.Range("A2:E30").ClearContents()

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
